$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.004257666666666667
$ws.Range("H2").Value = 0.012773
$ws.Range("M2").Value = 24.91851366666667
$ws.Range("N2").Value = 74.75554099999999
$ws.Range("O2").Value = 0.2924799159147552
$ws.Range("P2").Value = 0.2924799159147553
$ws.Range("Q2").Value = 0.1060947250214444
$ws.Range("R2").Value = 0.9548525251929999
$ws.Range("S2").Value = 0.2924799159147552
$ws.Range("T2").Value = 0.2924799159147553

# Row 3
$ws.Range("G3").Value = 0.004257666666666667
$ws.Range("H3").Value = 0.012773
$ws.Range("O3").Value = 0.4753125595076708
$ws.Range("P3").Value = 0.4753125595076708
$ws.Range("Q3").Value = 0.1724157884225556
$ws.Range("R3").Value = 1.551742095803
$ws.Range("S3").Value = 0.4753125595076708
$ws.Range("T3").Value = 0.4753125595076708

# Row 4
$ws.Range("G4").Value = 0.004257666666666667
$ws.Range("H4").Value = 0.012773
$ws.Range("M4").Value = 19.78346566666667
$ws.Range("N4").Value = 59.350397
$ws.Range("O4").Value = 0.232207524577574
$ws.Range("P4").Value = 0.232207524577574
$ws.Range("Q4").Value = 0.08423140232011112
$ws.Range("R4").Value = 0.758082620881
$ws.Range("S4").Value = 0.232207524577574
$ws.Range("T4").Value = 0.232207524577574
